$wb = $excel.ActiveWorkbook

# --- createUser sheet: bump the test-user numeric suffix 1071 -> 1072 ---
# (this also recalculates the dependent CONCAT formulas in B2 and F2)
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 1072
$wsCreateUser.Range("C17").Select() | Out-Null

# --- addListItem sheet: rename the list item UsersixteenB -> UsersixteenC ---
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "UsersixteenC"

# addListItem becomes the active sheet/tab, with A2 selected
$wsAddListItem.Activate() | Out-Null
$wsAddListItem.Range("A2").Select() | Out-Null
